$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 6 (ano_obj = 2025) metrics: total_customers, new_customers,
# new_rate and returning_rate reflect one additional new customer.
$ws.Range("C6").Value = 357
$ws.Range("E6").Value = 74
$ws.Range("G6").Value = 20.72829131652661
$ws.Range("H6").Value = 79.27170868347339
